$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row: value goes from empty to "false"
$ws.Range("B7").Value = "false"

# "Date" row: update timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"
